$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hashcode.csv")

$ws.Range("B9").Value = "b38f934c02d047a2ada11101a82c1f39"
$ws.Range("B17").Value = "07256692167359f375548b4159378639"
$ws.Range("B94").Value = "44213aeeab26b84a909d27da8747f1dd"
$ws.Range("B98").Value = "7e28e709da59e3fc566edfc13a487028"
$ws.Range("B109").Value = "4eadddab98df18409f53e51a7d916afb"
$ws.Range("B115").Value = "78fb34603fc974bb8815be6ff28d67f3"
$ws.Range("B159").Value = "7efd4d5ecec095ae0b2a2e3bc16c6c20"
$ws.Range("B169").Value = "413a0d05a619a60b898eb259c960afc0"
$ws.Range("B183").Value = "b3a5b41de62bc70708855999dc05272a"
$ws.Range("B200").Value = "167b2fa8a52251f81750b9c2cb5d4eea"
$ws.Range("B228").Value = "64b0b49079d4fafbf463562b0ce5c243"
$ws.Range("B293").Value = "9b7cf1e5faefbbd76ff6d5e598ea888b"
$ws.Range("B339").Value = "4dd4c1f8cdc1fd5cc6e0107860789455"
$ws.Range("B420").Value = "bf3569543f5afe0bd329968445d710df"
$ws.Range("B464").Value = "f3ca3a5e106381f567089cfeb1ff5eaa"
$ws.Range("B506").Value = "32b0e69ac96ccda0211b74f7e415d067"
$ws.Range("B524").Value = "e0be8f01f61a7e46740ea82661e2c46f"
$ws.Range("B580").Value = "fa0233183a94dd823d1a0c00a9af25d2"
$ws.Range("B600").Value = "98a7a4c7e45a4c7f13b04e8c8f695464"
$ws.Range("B604").Value = "ed159e1d6c2aca808fd7aad327e35968"
$ws.Range("B624").Value = "23a05fa1b6ac27eb97b8412b67c6f222"
$ws.Range("B635").Value = "31d4b27f68ee3e27be775bef84187400"
$ws.Range("B723").Value = "356ca7a6a0143f6e4c614d0549b08df8"
$ws.Range("B827").Value = "4cd16c911c9d83985478f327f616afa4"
$ws.Range("B882").Value = "c9c849f03081bb7a17b5eba5feebb7ea"
